$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9371385197543067
$ws.Range("J2").Value = 0.9371385197543067
$ws.Range("M2").Value = 110.642708
$ws.Range("N2").Value = 331.928124
$ws.Range("O2").Value = 0.5476418925386564
$ws.Range("P2").Value = 0.5476418925386564
$ws.Range("Q2").Value = 1526.492890145578
$ws.Range("R2").Value = 13738.43601131021
$ws.Range("S2").Value = 0.5132163125291236
$ws.Range("T2").Value = 0.5132163125291236

# Row 3
$ws.Range("I3").Value = 0.9371385197543067
$ws.Range("J3").Value = 0.9371385197543067
$ws.Range("O3").Value = 0.3151072754333865
$ws.Range("P3").Value = 0.3151072754333865
$ws.Range("S3").Value = 0.2952991656634565
$ws.Range("T3").Value = 0.2952991656634565

# Row 4
$ws.Range("I4").Value = 0.9371385197543067
$ws.Range("J4").Value = 0.9371385197543067
$ws.Range("M4").Value = 27.72944133333333
$ws.Range("N4").Value = 83.18832399999999
$ws.Range("O4").Value = 0.1372508320279571
$ws.Range("P4").Value = 0.1372508320279571
$ws.Range("Q4").Value = 382.5719363542897
$ws.Range("R4").Value = 3443.147427188608
$ws.Range("S4").Value = 0.1286230415617267
$ws.Range("T4").Value = 0.1286230415617267

# Row 5
$ws.Range("G5").Value = 0.9254496666666667
$ws.Range("H5").Value = 2.776349
$ws.Range("I5").Value = 0.0628614802456932
$ws.Range("J5").Value = 0.06286148024569319
$ws.Range("M5").Value = 110.642708
$ws.Range("N5").Value = 331.928124
$ws.Range("O5").Value = 0.5476418925386564
$ws.Range("P5").Value = 0.5476418925386564
$ws.Range("Q5").Value = 102.3942572376973
$ws.Range("R5").Value = 921.5483151392759
$ws.Range("S5").Value = 0.03442558000953278
$ws.Range("T5").Value = 0.03442558000953277

# Row 6
$ws.Range("G6").Value = 0.9254496666666667
$ws.Range("H6").Value = 2.776349
$ws.Range("I6").Value = 0.0628614802456932
$ws.Range("J6").Value = 0.06286148024569319
$ws.Range("O6").Value = 0.3151072754333865
$ws.Range("P6").Value = 0.3151072754333865
$ws.Range("Q6").Value = 58.91655817020723
$ws.Range("R6").Value = 530.2490235318651
$ws.Range("S6").Value = 0.01980810976993003
$ws.Range("T6").Value = 0.01980810976993003

# Row 7
$ws.Range("G7").Value = 0.9254496666666667
$ws.Range("H7").Value = 2.776349
$ws.Range("I7").Value = 0.0628614802456932
$ws.Range("J7").Value = 0.06286148024569319
$ws.Range("M7").Value = 27.72944133333333
$ws.Range("N7").Value = 83.18832399999999
$ws.Range("O7").Value = 0.1372508320279571
$ws.Range("P7").Value = 0.1372508320279571
$ws.Range("Q7").Value = 25.66220223878622
$ws.Range("R7").Value = 230.959820149076
$ws.Range("S7").Value = 0.008627790466230382
$ws.Range("T7").Value = 0.00862779046623038
